$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.103.07"
$ws.Range("E2").Value = "  -1.69%  "

$ws.Range("D3").Value = "1.895.72"
$ws.Range("E3").Value = "  -0.88%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'314.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").Value = "'0.5025"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.62%  "

$ws.Range("D8").Value = "'0.3906"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.18%  "

$ws.Range("D9").Value = "'0.09238"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.70%  "

$ws.Range("D10").Value = "'1.129"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.72%  "

$ws.Range("D11").Value = "'41.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.76%  "

$ws.Range("D12").Value = "'6.385"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.57%  "

$ws.Range("D13").Value = "'20.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.87%  "

$ws.Range("D14").Value = "1.904.30"
$ws.Range("E14").Value = "  -0.86%  "

$ws.Range("D15").Value = "'7.293"
$ws.Range("D15").Style = "Normal"

$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("D17").Value = "'92.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.65%  "

$ws.Range("D18").Value = "'0.00001108"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.84%  "

$ws.Range("D19").Value = "'0.06654"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").Value = "'17.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.64%  "

$ws.Range("D22").Value = "'6.213"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "28.162.75"
$ws.Range("E23").Value = "  -1.63%  "

$ws.Range("D24").Value = "'11.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("E25").Value = "  +1.91%  "

$ws.Range("D26").Value = "2.117.28"
$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("D27").Value = "'2.553"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.96%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.02%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'158.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("D30").Value = "'127.09"
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").Value = "'1.079"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.24%  "

$ws.Range("D32").Value = "'0.1057"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.47%  "

$ws.Range("D33").Value = "'5.613"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.50%  "

$ws.Range("D34").Value = "'3.609"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("D35").Value = "'9.547"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.19%  "

$ws.Range("D36").Value = "'1.370"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +15.43%  "

$ws.Range("D37").Value = "'0.06607"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.81%  "

$ws.Range("D38").Value = "'0.02403"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.65%  "

$ws.Range("D39").Value = "'0.2205"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.35%  "

$ws.Range("D40").Value = "'1.223"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.01%  "

$ws.Range("D41").Value = "'0.6468"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.37%  "

$ws.Range("D42").Value = "'11.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.72%  "

$ws.Range("D43").Value = "'4.976"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.45%  "

$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").Value = "'0.6090"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").Value = "'13.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.68%  "

$ws.Range("D47").Value = "'1.302"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.54%  "

$ws.Range("D48").Value = "'3.692"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.62%  "

$ws.Range("D49").Value = "'2.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.09%  "

$ws.Range("D50").Value = "'122.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.33%  "

$ws.Range("D51").Value = "'1.199"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.26%  "
